# data saving for individual level
#
# Updates the "enddate" column (F) for two rows:
#   - F2: "2024-08-03" -> "2024-08-10"
#   - F9: (empty)      -> "2024-08-10"
#
# Both cells are formatted with a date number format (style index 1).
# Assigning a date-looking string straight to `.Value` on a date-formatted
# cell causes Excel to parse it into a date serial number instead of
# keeping the literal text, and forcing text (via NumberFormat "@" or a
# leading apostrophe) mints a brand new style entry, which would change
# the cell's style index. To preserve both the literal text value *and*
# the original style (s="1"), we stage the literal text in a scratch
# cell (forcing it to text with a leading apostrophe), copy it, and use
# PasteSpecial (values only) onto the target cell - this pastes only the
# value and leaves the destination's existing formatting untouched. The
# scratch cell is restored to its original content afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

function Set-LiteralText($range, [string]$text) {
    $scratch = $ws.Range("H17")
    $savedValue = $scratch.Value2

    $scratch.Value = "'" + $text
    $scratch.Copy()
    $range.PasteSpecial($xlPasteValues)

    $scratch.Value = $savedValue
}

Set-LiteralText $ws.Range("F2") "2024-08-10"
Set-LiteralText $ws.Range("F9") "2024-08-10"
